$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.185.01"
$ws.Range("E2").Value = "  +2.59%  "
$ws.Range("D3").Value = "3.449.03"
$ws.Range("E3").Value = "  +1.71%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "'578.62"
$ws.Range("E5").Value = "  +3.49%  "
$ws.Range("D6").Value = "'187.65"
$ws.Range("E6").Value = "  +6.63%  "
$ws.Range("D7").Value = "'0.630"
$ws.Range("E7").Value = "  +0.05%  "
$ws.Range("D8").Value = "3.443.19"
$ws.Range("E8").Value = "  +1.79%  "
$ws.Range("D9").Value = "'1.00"
$ws.Range("E9").Value = "  -0.07%  "
$ws.Range("E10").Value = "  +0.59%  "
$ws.Range("E11").Value = "  +0.52%  "
$ws.Range("D12").Value = "'57.64"
$ws.Range("E12").Value = "  +7.87%  "
$ws.Range("E13").Value = "  -0.90%  "
$ws.Range("D14").Value = "'9.43"
$ws.Range("E14").Value = "  +2.23%  "
$ws.Range("D15").Value = "3.994.58"
$ws.Range("E15").Value = "  +1.55%  "
$ws.Range("D16").Value = "'19.01"
$ws.Range("E16").Value = "  +3.83%  "
$ws.Range("D17").Value = "3.440.28"
$ws.Range("E17").Value = "  +1.34%  "
$ws.Range("D18").Value = "67.102.54"
$ws.Range("E18").Value = "  +2.55%  "
$ws.Range("E19").Value = "  -0.61%  "
$ws.Range("D20").Value = "'12.05"
$ws.Range("E20").Value = "  +1.72%  "
$ws.Range("E21").Value = "  +1.60%  "
$ws.Range("D22").Value = "'489.32"
$ws.Range("E22").Value = "  +4.21%  "
$ws.Range("D23").Value = "'5.65"
$ws.Range("E23").Value = "  +14.12%  "
$ws.Range("D24").Value = "'17.02"
$ws.Range("E24").Value = "  +18.93%  "
$ws.Range("D25").Value = "'4.33"
$ws.Range("E25").Value = "  +4.96%  "
$ws.Range("D26").Value = "'89.54"
$ws.Range("E26").Value = "  +2.60%  "
$ws.Range("D27").Value = "'2.97"
$ws.Range("E27").Value = "  +1.88%  "
$ws.Range("D28").Value = "'10.94"
$ws.Range("E28").Value = "  +2.05%  "
$ws.Range("E29").Value = "  +3.17%  "
$ws.Range("D30").Value = "'31.21"
$ws.Range("E30").Value = "  +0.45%  "
$ws.Range("D31").Value = "'7.33"
$ws.Range("E31").Value = "  +11.89%  "
$ws.Range("D32").Value = "'604.89"
$ws.Range("E32").Value = "  +5.46%  "
$ws.Range("D33").Value = "'64.73"
$ws.Range("E33").Value = "  +2.58%  "
$ws.Range("D34").Value = "'11.79"
$ws.Range("E34").Value = "  +2.57%  "
$ws.Range("E35").Value = "  +3.36%  "
$ws.Range("E36").Value = "  +0.00%  "
$ws.Range("E37").Value = "  +3.24%  "
$ws.Range("D38").Value = "'37.01"
$ws.Range("E38").Value = "  +2.92%  "
$ws.Range("D39").Value = "0.0₃0781"
$ws.Range("E39").Value = "  +5.47%  "
$ws.Range("E40").Value = "  +3.23%  "
$ws.Range("E41").Value = "  -4.21%  "
$ws.Range("D42").Value = "3.190.51"
$ws.Range("E42").Value = "  +2.14%  "
$ws.Range("E43").Value = "  +3.19%  "
$ws.Range("E44").Value = "  +2.84%  "
$ws.Range("D45").Value = "'2.57"
$ws.Range("E45").Value = "  +5.42%  "
$ws.Range("D46").Value = "'3.22"
$ws.Range("E46").Value = "  +1.70%  "
$ws.Range("E47").Value = "  +1.23%  "
$ws.Range("D48").Value = "'2.66"
$ws.Range("E48").Value = "  +14.98%  "
$ws.Range("D49").Value = "'0.999"
$ws.Range("E49").Value = "  -0.02%  "
$ws.Range("D50").Value = "'8.62"
$ws.Range("E50").Value = "  +2.32%  "
$ws.Range("D51").Value = "'140.48"
$ws.Range("E51").Value = "  -0.09%  "
